# Slide 3 ("Aula 11 - Algoritmos e Complexidade - Divisão e Conquista"),
# body placeholder: bold the word "final" in
#   "...combinar as soluções para obter a solução final."
# and add a trailing space after the final period, splitting the
# trailing ". " off into its own run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Bold just the word "final" (this splits the run it lives in).
$finalWord = $tr.Find("final")
$finalWord.Font.Bold = $true

# Re-acquire the range (text length is unchanged) and grab the trailing
# period, then add a trailing space after it.
$tr2 = $sh.TextFrame.TextRange
$period = $tr2.Characters($tr2.Length, 1)
$null = $period.InsertAfter(" ")
